# 9.5.2.xlsx update: add a new "2021" / "515" data column (O) mirroring the
# existing "2020" / "534" column (N), and move the sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new year column (O) -------------------------------------------
# Row 4 holds the year headers, row 5 holds the corresponding values.
$ws.Cells.Item(4, 15).Value = 2021
$ws.Cells.Item(5, 15).Value = 515

# Copy the formatting from the preceding column (N) so the new cells pick up
# the same styles (right-aligned year header / numeric value cell).
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)   # xlPasteFormats

# --- Update the view/selection ----------------------------------------------
# Select P12 so it becomes the active cell/selection and the view scrolls
# back to show column A (clearing the previous topLeftCell="E1").
$ws.Range("P12").Select() | Out-Null
